# Add Sum comparison to MasterGallery sheets
$wb = $excel.ActiveWorkbook

# --- "Formula Samples" sheet: add a Sum row (row 11) ---
$wsFormula = $wb.Worksheets.Item("Formula Samples")
$wsFormula.Range("A11").Value = "Sum"
$wsFormula.Range("B11").Formula = "=SUM(B2:B4)"

# Add the rubric comment describing how B11 should be graded.
$sumComment = $wsFormula.Range("B11").AddComment("rubric:`n score: 2`n type: formula")

# Update the selection on "Formula Samples" (no longer the active tab).
[void]$wsFormula.Range("B12").Select()

# --- "Formula Samples_CheckOrder" sheet: register the new B11 cell ---
$wsFormulaCheck = $wb.Worksheets.Item("Formula Samples_CheckOrder")
$wsFormulaCheck.Range("A8").Value = 7
$wsFormulaCheck.Range("B8").Value = "B11"

# Selecting here makes "Formula Samples_CheckOrder" the active tab,
# matching the workbook's new activeTab, and clears the previous
# "Test Case Samples_CheckOrder" tab selection.
[void]$wsFormulaCheck.Range("A9").Select()
